$wb = $excel.ActiveWorkbook

# Rename the worksheets (Лист1/2/3 -> Worksheet 1/2/3)
$wb.Worksheets.Item(1).Name = "Worksheet 1"
$wb.Worksheets.Item(2).Name = "Worksheet 2"
$wb.Worksheets.Item(3).Name = "Worksheet 3"

# Update the date text on the first worksheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Date: 05-10-2018 - Department: Sales department"
